$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = "44d95e55ae00b14a3d720888e568e870"
$ws.Range("B49").Value = "ec239f2beb1fd4c210927b736b7d41d0"
$ws.Range("B65").Value = "0159658a61ec1e54769b3e111cd5e338"
$ws.Range("B80").Value = "5c7461dca3df71bee93b2ffb4c7aff6b"
$ws.Range("B81").Value = "bc13665eac829680b6a0efac910209a9"
$ws.Range("B113").Value = "d2c1e61c41dbe7d42161673c62f26d94"
$ws.Range("B117").Value = "c44933a8687ca715bd1e53da6d63de28"
$ws.Range("B122").Value = "6403db4eaca423e88668dcf035f15b05"
$ws.Range("B136").Value = "87db705baeeef194c934e2da535c0206"
$ws.Range("B163").Value = "e6bcb2e1d0134c14259cc37457c02026"
$ws.Range("B187").Value = "bbe42b101b0df60ce51714a09105540e"
$ws.Range("B527").Value = "bbe42b101b0df60ce51714a09105540e"
$ws.Range("B191").Value = "78eaa3980b03470a955812e5188ad493"
$ws.Range("B198").Value = "10f65bff40c998d991b943c73cd44ddd"
$ws.Range("B222").Value = "8ba01cccfb66474d42056f6f56f8c832"
$ws.Range("B232").Value = "d53ce3f95a81937bf61b93ea482d1bb9"
$ws.Range("B296").Value = "3e7d4df85204ac4351dfbec350fa0213"
$ws.Range("B342").Value = "091253ce18fb2eaae2c3e52d7191a868"
$ws.Range("B350").Value = "205045de71ccf4d8ebb7043be63d7d1e"
$ws.Range("B360").Value = "8378c8ce3a4390b4106ae67049b24cbb"
$ws.Range("B404").Value = "efc28c664155744781a804fccc281880"
$ws.Range("B419").Value = "bf3569543f5afe0bd329968445d710df"
$ws.Range("B465").Value = "3fb2c89a1ead2341473dd5438758784d"
$ws.Range("B480").Value = "5b827e9d111ec4a61c8f3b1afd227402"
$ws.Range("B487").Value = "98c3b60d2a18c880caa013fd2a9fb35b"
$ws.Range("B507").Value = "caf0902acd5e4ab007abd4dbb31c1a66"
$ws.Range("B514").Value = "e50e7ab7e04d7767992587dc16d148d2"
$ws.Range("B548").Value = "eeae9bdd031c807456d6814fc03e4b4e"
$ws.Range("B574").Value = "f3ec5be5d08573163925e4f336c4139c"
$ws.Range("B619").Value = "5f08ecffecd63e81c0870b802b54b76b"
$ws.Range("B697").Value = "5d36b099269766af24b39475fff3e8e3"
$ws.Range("B700").Value = "82b72890bcd18c3586e86b913781f104"
$ws.Range("B726").Value = "9bc32158d70334da0475337d0e357327"
$ws.Range("B763").Value = "b8336d370aee1c9076bd82993b452960"
$ws.Range("B822").Value = "d5b471cb70b49bc80a9c2c37ae03f627"
$ws.Range("B854").Value = "523e4ae528e99719f3dd21cbb14652cd"
$ws.Range("B889").Value = "ec565bb99879f865a731ab258df28300"
$ws.Range("B890").Value = "e731922839d98123f930456d64d6c8d3"
$ws.Range("B946").Value = "164564ca6182282ff0c3c6b63f6c25c6"
